$d = $word.ActiveDocument

# The existing "Compact" paragraph style is being duplicated into a new
# "CompactList" style (display name "Compact List") with identical
# formatting: based on Body Text, quick style, 36-twip (1.8pt) spacing
# before/after.
$source = $d.Styles.Item("Compact")

$newStyle = $d.Styles.Add("Compact List", $source.Type)   # wdStyleTypeParagraph

# NOTE: assign the base style by its style id (not a Style object / its
# NameLocal) so it serializes as <w:basedOn w:val="BodyText"/> rather than
# the display name "Body Text".
$newStyle.BaseStyle = "BodyText"
$newStyle.QuickStyle = $source.QuickStyle
$newStyle.ParagraphFormat.SpaceBefore = $source.ParagraphFormat.SpaceBefore
$newStyle.ParagraphFormat.SpaceAfter = $source.ParagraphFormat.SpaceAfter
